# feat: add 2022-Q1 data
#
# The workbook has five quarterly sheets ending in a "总计" (totals) summary
# sheet. This adds a new "2022-Q1" sheet (in the same per-fund layout as the
# other quarter sheets) right before "总计", and updates "总计" with a new
# leading row summarizing the 2022-Q1 quarter.

$wb = $excel.ActiveWorkbook

# --- Step 1: turn the existing "总计" sheet into the new "2022-Q1" sheet ---
# Reusing this worksheet (instead of inserting a brand new one) keeps the
# existing header/column styling (bold, bordered, centered - style index 2)
# intact for the cells that are shared with the other quarter sheets.
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# The old header only spanned B1:D1 - copy that cell's format onto the newly
# used E1:H1 header cells so they pick up the same style (index 2).
$q1.Range("B1").Copy() | Out-Null
$q1.Range("E1:H1").PasteSpecial(-4122) | Out-Null

# Fund holdings data rows (A is the existing styled 0-based index column,
# H is the unstyled rank column - both numeric; B/C/D/E/F/G are text)
$fundRows = @(
    @(0, "009837", "华夏磐锐一年定期开放混合A",       "16.45", "79.44", "3.71", "0.6103", 4),
    @(1, "005947", "德邦民裕进取量化精选灵活配置混合A", "0.53",  "94.44", "7.29", "0.0386", 3),
    @(2, "009838", "华夏磐锐一年定期开放混合C",       "0.44",  "79.44", "3.71", "0.0163", 4),
    @(3, "014331", "华泰柏瑞中证稀土产业ETF联接A",    "0.86",  "24.22", "1.00", "0.0086", 10),
    @(4, "014332", "华泰柏瑞中证稀土产业ETF联接C",    "0.70",  "24.22", "1.00", "0.0070", 10),
    @(5, "005948", "德邦民裕进取量化精选灵活配置混合C", "0.09",  "94.44", "7.29", "0.0066", 3),
    @(6, "002020", "国都创新驱动灵活配置混合",        "0.15",  "74.87", "2.55", "0.0038", 5)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    $q1.Range("A$r").Value = $row[0]

    # B..G are stored as text (fund codes must keep leading zeros, and the
    # numeric-looking figures are text in the source data too) - the leading
    # apostrophe forces Excel to store them as text instead of numbers.
    $q1.Range("B$r").Value = "'" + $row[1]
    $q1.Range("C$r").Value = "'" + $row[2]
    $q1.Range("D$r").Value = "'" + $row[3]
    $q1.Range("E$r").Value = "'" + $row[4]
    $q1.Range("F$r").Value = "'" + $row[5]
    $q1.Range("G$r").Value = "'" + $row[6]

    $q1.Range("H$r").Value = $row[7]
}

# --- Step 2: add a brand-new "总计" sheet after "2022-Q1" -----------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @("2022-Q1", 7, 0.6899999999999999),
    @("2021-Q4", 4, 0.12),
    @("2021-Q3", 2, 0.13),
    @("2021-Q1", 4, 0.03),
    @("2020-Q4", 2, 0.01)
)

for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]
    $total.Range("A$r").Value = $i
    $total.Range("B$r").Value = $row[0]
    $total.Range("C$r").Value = $row[1]
    $total.Range("D$r").Value = $row[2]
}

# Match the original style (index 2: bold/bordered/centered) used on the
# header row and column A throughout the workbook's quarter sheets.
$q1.Range("A2").Copy() | Out-Null
$total.Range("A2:A6").PasteSpecial(-4122) | Out-Null
$q1.Range("B1").Copy() | Out-Null
$total.Range("B1:D1").PasteSpecial(-4122) | Out-Null

# --- Step 3: restore the originally-active sheet ---------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
